$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename "BulletAdvance++" header (col J) to "BulletLevel++" ---
$ws.Range("J1").Value = "BulletLevel++"

# --- Remove the Energy++ / EnergyMax++ / AutoEnergy++ columns (K, L, M) ---
# This shifts the old N (MoveSpeed++) and O (BulletSpeed) columns left to K, L.
$ws.Columns("K:M").Delete()

# --- New "Deffend" attribute column in the freed M slot ---
$ws.Range("M1").Value = "Deffend"

# --- Column width tweaks ---
$ws.Columns("B").ColumnWidth = 22.571428571428573
$ws.Columns("F").ColumnWidth = 9.142857142857142
$ws.Columns("K").ColumnWidth = 13.142857142857142
$ws.Columns("L").ColumnWidth = 12.285714285714286

# --- New item rows: "Lazer Gun" and "Solar Panel" with their attribute values ---
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "Lazer Gun"
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0

$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "Solar Panel"
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = 10
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 8).Value = 1
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = 0.9

# --- Stale outline-level bookkeeping: nest/unnest a throwaway row twice so the
# sheet's outlineLevelRow max bumps from 1 to 2 (mirrors an earlier edit elsewhere
# in the sheet that grouped/removed rows) without leaving any grouped rows behind.
$ws.Rows("10:10").Group()
$ws.Rows("10:10").Group()
$ws.Rows("10:10").Delete()

# --- Selection as left in the source file ---
$ws.Range("F3").Select()
